$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 14.956293
$ws.Range("H2").Value = 44.868879
$ws.Range("I2").Value = 0.195346866793292
$ws.Range("J2").Value = 0.1953468667932919
$ws.Range("M2").Value = 2.318119
$ws.Range("N2").Value = 6.954357
$ws.Range("O2").Value = 0.1070970465647729
$ws.Range("P2").Value = 0.1070970465647729
$ws.Range("Q2").Value = 34.670466972867
$ws.Range("R2").Value = 312.034202755803
$ws.Range("S2").Value = 0.02092107248924368
$ws.Range("T2").Value = 0.02092107248924368

$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 14.956293
$ws.Range("H3").Value = 44.868879
$ws.Range("I3").Value = 0.195346866793292
$ws.Range("J3").Value = 0.1953468667932919
$ws.Range("N3").Value = 5.768654000000001
$ws.Range("O3").Value = 0.0888372291002696
$ws.Range("P3").Value = 0.0888372291002696
$ws.Range("Q3").Value = 28.759226479874
$ws.Range("R3").Value = 258.833038318866
$ws.Range("S3").Value = 0.01735407435933553
$ws.Range("T3").Value = 0.01735407435933552

$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 14.956293
$ws.Range("H4").Value = 44.868879
$ws.Range("I4").Value = 0.195346866793292
$ws.Range("J4").Value = 0.1953468667932919
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 0.6666666666666666
$ws.Range("M4").Value = 0.108471
$ws.Range("N4").Value = 0.325413
$ws.Range("O4").Value = 0.005011357802566427
$ws.Range("P4").Value = 0.005011357802566428
$ws.Range("Q4").Value = 1.622324058003
$ws.Range("R4").Value = 14.600916522027
$ws.Range("S4").Value = 0.000978953045111468
$ws.Range("T4").Value = 0.000978953045111468

$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 14.956293
$ws.Range("H5").Value = 44.868879
$ws.Range("I5").Value = 0.195346866793292
$ws.Range("J5").Value = 0.1953468667932919
$ws.Range("M5").Value = 17.29555733333333
$ws.Range("N5").Value = 51.886672
$ws.Range("O5").Value = 0.7990543665323911
$ws.Range("P5").Value = 0.7990543665323911
$ws.Range("Q5").Value = 258.677423075632
$ws.Range("R5").Value = 2328.096807680688
$ws.Range("S5").Value = 0.1560927668996013
$ws.Range("T5").Value = 0.1560927668996013

$ws.Range("I6").Value = 0.653630228225219
$ws.Range("J6").Value = 0.653630228225219
$ws.Range("M6").Value = 2.318119
$ws.Range("N6").Value = 6.954357
$ws.Range("O6").Value = 0.1070970465647729
$ws.Range("P6").Value = 0.1070970465647729
$ws.Range("Q6").Value = 116.0073136168067
$ws.Range("R6").Value = 1044.06582255126
$ws.Range("S6").Value = 0.07000186698837942
$ws.Range("T6").Value = 0.07000186698837943

$ws.Range("I7").Value = 0.653630228225219
$ws.Range("J7").Value = 0.653630228225219
$ws.Range("N7").Value = 5.768654000000001
$ws.Range("O7").Value = 0.0888372291002696
$ws.Range("P7").Value = 0.0888372291002696
$ws.Range("Q7").Value = 96.22831467019111
$ws.Range("R7").Value = 866.0548320317201
$ws.Range("S7").Value = 0.05806669833170529
$ws.Range("T7").Value = 0.05806669833170529

$ws.Range("I8").Value = 0.653630228225219
$ws.Range("J8").Value = 0.653630228225219
$ws.Range("K8").Value = 2
$ws.Range("L8").Value = 0.6666666666666666
$ws.Range("M8").Value = 0.108471
$ws.Range("N8").Value = 0.325413
$ws.Range("O8").Value = 0.005011357802566427
$ws.Range("P8").Value = 0.005011357802566428
$ws.Range("Q8").Value = 5.42829307526
$ws.Range("R8").Value = 48.85463767734
$ws.Range("S8").Value = 0.003275574944209725
$ws.Range("T8").Value = 0.003275574944209726

$ws.Range("I9").Value = 0.653630228225219
$ws.Range("J9").Value = 0.653630228225219
$ws.Range("M9").Value = 17.29555733333333
$ws.Range("N9").Value = 51.886672
$ws.Range("O9").Value = 0.7990543665323911
$ws.Range("P9").Value = 0.7990543665323911
$ws.Range("Q9").Value = 865.5341437369956
$ws.Range("R9").Value = 7789.807293632961
$ws.Range("S9").Value = 0.5222860879609246
$ws.Range("T9").Value = 0.5222860879609246

$ws.Range("G10").Value = 10.01531766666667
$ws.Range("H10").Value = 30.045953
$ws.Range("I10").Value = 0.1308118880876991
$ws.Range("J10").Value = 0.1308118880876991
$ws.Range("M10").Value = 2.318119
$ws.Range("N10").Value = 6.954357
$ws.Range("O10").Value = 0.1070970465647729
$ws.Range("P10").Value = 0.1070970465647729
$ws.Range("Q10").Value = 23.21669817413567
$ws.Range("R10").Value = 208.950283567221
$ws.Range("S10").Value = 0.01400956686975417
$ws.Range("T10").Value = 0.01400956686975417

$ws.Range("G11").Value = 10.01531766666667
$ws.Range("H11").Value = 30.045953
$ws.Range("I11").Value = 0.1308118880876991
$ws.Range("J11").Value = 0.1308118880876991
$ws.Range("N11").Value = 5.768654000000001
$ws.Range("O11").Value = 0.0888372291002696
$ws.Range("P11").Value = 0.0888372291002696
$ws.Range("Q11").Value = 19.25830077302911
$ws.Range("R11").Value = 173.324706957262
$ws.Range("S11").Value = 0.01162096567108575
$ws.Range("T11").Value = 0.01162096567108575

$ws.Range("G12").Value = 10.01531766666667
$ws.Range("H12").Value = 30.045953
$ws.Range("I12").Value = 0.1308118880876991
$ws.Range("J12").Value = 0.1308118880876991
$ws.Range("K12").Value = 2
$ws.Range("L12").Value = 0.6666666666666666
$ws.Range("M12").Value = 0.108471
$ws.Range("N12").Value = 0.325413
$ws.Range("O12").Value = 0.005011357802566427
$ws.Range("P12").Value = 0.005011357802566428
$ws.Range("Q12").Value = 1.086371522621
$ws.Range("R12").Value = 9.777343703589
$ws.Range("S12").Value = 0.0006555451760367368
$ws.Range("T12").Value = 0.0006555451760367369

$ws.Range("G13").Value = 10.01531766666667
$ws.Range("H13").Value = 30.045953
$ws.Range("I13").Value = 0.1308118880876991
$ws.Range("J13").Value = 0.1308118880876991
$ws.Range("M13").Value = 17.29555733333333
$ws.Range("N13").Value = 51.886672
$ws.Range("O13").Value = 0.7990543665323911
$ws.Range("P13").Value = 0.7990543665323911
$ws.Range("Q13").Value = 173.2205009153796
$ws.Range("R13").Value = 1558.984508238416
$ws.Range("S13").Value = 0.1045258103708224
$ws.Range("T13").Value = 0.1045258103708224

$ws.Range("G14").Value = 1.547411
$ws.Range("H14").Value = 4.642232999999999
$ws.Range("I14").Value = 0.0202110168937901
$ws.Range("J14").Value = 0.0202110168937901
$ws.Range("M14").Value = 2.318119
$ws.Range("N14").Value = 6.954357
$ws.Range("O14").Value = 0.1070970465647729
$ws.Range("P14").Value = 0.1070970465647729
$ws.Range("Q14").Value = 3.587082839908999
$ws.Range("R14").Value = 32.28374555918099
$ws.Range("S14").Value = 0.00216454021739565
$ws.Range("T14").Value = 0.00216454021739565

$ws.Range("G15").Value = 1.547411
$ws.Range("H15").Value = 4.642232999999999
$ws.Range("I15").Value = 0.0202110168937901
$ws.Range("J15").Value = 0.0202110168937901
$ws.Range("N15").Value = 5.768654000000001
$ws.Range("O15").Value = 0.0888372291002696
$ws.Range("P15").Value = 0.0888372291002696
$ws.Range("Q15").Value = 2.975492884931333
$ws.Range("R15").Value = 26.779435964382
$ws.Range("S15").Value = 0.001795490738143051
$ws.Range("T15").Value = 0.00179549073814305

$ws.Range("G16").Value = 1.547411
$ws.Range("H16").Value = 4.642232999999999
$ws.Range("I16").Value = 0.0202110168937901
$ws.Range("J16").Value = 0.0202110168937901
$ws.Range("K16").Value = 2
$ws.Range("L16").Value = 0.6666666666666666
$ws.Range("M16").Value = 0.108471
$ws.Range("N16").Value = 0.325413
$ws.Range("O16").Value = 0.005011357802566427
$ws.Range("P16").Value = 0.005011357802566428
$ws.Range("Q16").Value = 0.167849218581
$ws.Range("R16").Value = 1.510642967229
$ws.Range("S16").Value = 0.0001012846372084969
$ws.Range("T16").Value = 0.0001012846372084969

$ws.Range("G17").Value = 1.547411
$ws.Range("H17").Value = 4.642232999999999
$ws.Range("I17").Value = 0.0202110168937901
$ws.Range("J17").Value = 0.0202110168937901
$ws.Range("M17").Value = 17.29555733333333
$ws.Range("N17").Value = 51.886672
$ws.Range("O17").Value = 0.7990543665323911
$ws.Range("P17").Value = 0.7990543665323911
$ws.Range("Q17").Value = 26.76333566873066
$ws.Range("R17").Value = 240.870021018576
$ws.Range("S17").Value = 0.0161497013010429
$ws.Range("T17").Value = 0.0161497013010429

